$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 and 3 (the two "N/A" placeholder rows), shifting remaining
# rows up so the real tender data starts at row 2.
$ws.Range("A2:B3").EntireRow.Delete()
